# Update the "Förändrad" (Changed) date column C for rows 2-29
# from 2024-04-14 (serial 45396) to 2024-04-18 (serial 45400).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 29; $row++) {
    $ws.Cells.Item($row, 3).Value = 45400
}
